$wb = $excel.ActiveWorkbook

$oldGuid = "25f6d518-2ea6-423c-b08c-74525fa1b5f1"
$newGuid = "d3e98c78-52c1-48c9-870d-bc0ad3910f2d"

$oldHoHash = "5e106ced215fb5dc13f6056d0c067291906a3dc5"
$newHoHash = "5844b2e5a7e0a2a2313eaf3fd2d3614df64b38a7"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-09-03 03:01:53"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHoHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-03 03:01:49"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHoHash.de-de.xlf"

# --- Update hyperlink display text (keep same target) ---
foreach ($ws in @($wsOverview, $wsZhCn, $wsDeDe)) {
    foreach ($hl in $ws.Hyperlinks) {
        $hl.TextToDisplay = $hl.TextToDisplay -replace [regex]::Escape($oldGuid), $newGuid
    }
}
